# editSyncOptions.xlsx - "All feedback changes done"
#
# Content changes:
#  - B1 header "Host Number" -> "Host Name"
#  - A3/A4 "One For All" -> "Second Flow"
#  - B3/B4 numeric host numbers (1/2) -> text source host names
#    ("psp-MyLinSecondFlow-src1" / "psp-MyLinSecondFlow-src2")
#  - Column B / S best-fit widths change as a side effect of the wider
#    text now stored in column B (and a minor manual resize of column S)
#  - Selection moves to T4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new source-host values first (B3/B4) so the new shared strings
# are appended to sharedStrings.xml in the same order the source
# workbook uses (src1, src2, then the Host Name / Second Flow labels).
$ws.Range("B3").Value = "psp-MyLinSecondFlow-src1"
$ws.Range("B4").Value = "psp-MyLinSecondFlow-src2"

# Rename the "Host Number" column header to "Host Name".
$ws.Range("B1").Value = "Host Name"

# Rename the wave from "One For All" to "Second Flow".
$ws.Range("A3").Value = "Second Flow"
$ws.Range("A4").Value = "Second Flow"

# Column B grew to fit the longer host names; column S was nudged a
# touch narrower. (ColumnWidth is quantized to the workbook's default
# font metrics, so these are the closest attainable character widths.)
$ws.Columns.Item(2).ColumnWidth = 22.666666666666664
$ws.Columns.Item(19).ColumnWidth = 12.333333333333332

# Selection ends on T4.
$ws.Range("T4").Select() | Out-Null
